$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- semana 5 de 2026 -------------------------------------------------
# A new record ("620" Parotiditis) was reported for this week, so a row
# is inserted at position 25; the former rows 25-29 (Sifilis gestacional
# ... Vih/sida) shift down to rows 26-30.
$ws.Rows.Item(25).Insert()

# New row: 620 / Parotiditis
$ws.Range("A25").Value = "'620"
$ws.Range("B25").Value = "Parotiditis"
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0.27

# Refresh Esperado (C), Observado (D) and valor p (E) for every event
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.15
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 0.09
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0.18
$ws.Range("C7").Value = 47
$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 0.04
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.37
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 0.06
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 0.01
$ws.Range("C11").Value = 154
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 0.04
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.05
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0.04
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0.37
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 1
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 1
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 0.07000000000000001
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.37
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 1
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0.22
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 1
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 0.05
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 0.09
$ws.Range("C30").Value = 9
$ws.Range("D30").Value = 8
$ws.Range("E30").Value = 0.13
